$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Helper: give a cell a "top+bottom only" thin border (reuses the
# existing borderId=4 definition in styles.xml) starting from a
# clean slate (default font, no alignment) rather than inheriting
# the bold/centered look of the cell's previous style.
#
# The whole-collection ".Borders.Weight = " assignment commits a
# single all-sides-thin style in one shot (matching the existing
# borderId=1), after which clearing the unwanted edges one at a
# time still only ever passes through already-existing border
# geometries (no new/orphan style or border entries get created).
# ------------------------------------------------------------------
function Set-TopBottomBorder($cell) {
    $cell.ClearFormats()
    $cell.Borders.Weight = 2              # all edges thin (borderId=1)
    $cell.Borders.Item(7).LineStyle = 0   # xlEdgeLeft  -> none
    $cell.Borders.Item(10).LineStyle = 0  # xlEdgeRight -> none
}

# ------------------------------------------------------------------
# Helper: give a cell a "top+right+bottom only" thin border (reuses
# the existing borderId=5 definition), same clean-slate approach.
# ------------------------------------------------------------------
function Set-TopRightBottomBorder($cell) {
    $cell.ClearFormats()
    $cell.Borders.Weight = 2              # all edges thin (borderId=1)
    $cell.Borders.Item(7).LineStyle = 0   # xlEdgeLeft -> none
}

# ------------------------------------------------------------------
# Sheet 1: quality_comparison
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-TopBottomBorder      $ws1.Range("C1")
Set-TopRightBottomBorder $ws1.Range("D1")

$ws1.Range("C2").Value = "approach"

# ------------------------------------------------------------------
# Sheet 2: computational_comparison
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-TopBottomBorder      $ws2.Range("C1")
Set-TopRightBottomBorder $ws2.Range("D1")
Set-TopBottomBorder      $ws2.Range("F1")
Set-TopRightBottomBorder $ws2.Range("G1")

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell at G5.
$ws2.Range("G5").ClearContents()
